# Weekly fruit/vegetable price update:
# Insert a new data row at row 93 (pushing the existing rows 93-114 down to
# 94-115) and populate it with this week's record for Espinaca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 93:114 down to 94:115, carrying formatting (incl. the date
# number format on column D) along with them.
$ws.Rows("93:93").Insert()

# Populate the newly inserted row 93 with the new weekly record.
$ws.Range("A93").Value = 1
$ws.Range("B93").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C93").Value = "Arica y Parinacota"
$ws.Range("D93").Value = 45204
$ws.Range("E93").Value = 15
$ws.Range("F93").Value = 100112012
$ws.Range("G93").Value = "Espinaca"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Segunda"
$ws.Range("J93").Value = 250
$ws.Range("K93").Value = 1400
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 1450
$ws.Range("N93").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 483
$ws.Range("Q93").Value = 3
$ws.Range("R93").Value = "Hortaliza"
